$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-order the three "Criterion N, Air Speed 0.1" tabs.
#    Old tab order (positions 3,4,5): Criterion 2, Criterion 1, Criterion 3
#    New tab order (positions 3,4,5): Criterion 1, Criterion 3, Criterion 2
#    Using Worksheet.Move physically rotates each tab's row data along with
#    it, so the sheet that wiped up in slot 3 now carries what used to be
#    sheet 4's data (and so on).
# ---------------------------------------------------------------------------
$crit1 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$crit2 = $wb.Worksheets.Item("Criterion 2, Air Speed 0.1")
$crit1.Move($crit2)

$crit3 = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")
$crit2b = $wb.Worksheets.Item("Criterion 2, Air Speed 0.1")
$crit3.Move($crit2b)

# ---------------------------------------------------------------------------
# 2. Fix up the "Absolute/Relative Change" table headers on each of the
#    rotated sheets so the column names refer to the criterion that is now
#    actually shown on that tab.
# ---------------------------------------------------------------------------
$wsCrit1 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")
$wsCrit1.Cells.Item(1, 5).Value = "Criterion 1 Absolute Change"
$wsCrit1.Cells.Item(1, 6).Value = "Criterion 1 Relative Change (%)"

$wsCrit3 = $wb.Worksheets.Item("Criterion 3, Air Speed 0.1")
$wsCrit3.Cells.Item(1, 5).Value = "Criterion 3 Absolute Change"
$wsCrit3.Cells.Item(1, 6).Value = "Criterion 3 Relative Change (%)"

$wsCrit2 = $wb.Worksheets.Item("Criterion 2, Air Speed 0.1")
$wsCrit2.Cells.Item(1, 5).Value = "Criterion 2 Absolute Change"
$wsCrit2.Cells.Item(1, 6).Value = "Criterion 2 Relative Change (%)"

# ---------------------------------------------------------------------------
# 3. Update the "readme" log table: re-order its columns from
#    index,JobNo,Date,Author,sheet_name to index,Author,sheet_name,JobNo,Date
#    and refresh the data to match (new run date + the sheet_name values
#    reflecting the sheet tabs in their new order).
# ---------------------------------------------------------------------------
$wsReadme = $wb.Worksheets.Item("readme")

$wsReadme.Cells.Item(1, 2).Value = "Author"
$wsReadme.Cells.Item(1, 3).Value = "sheet_name"
$wsReadme.Cells.Item(1, 4).Value = "JobNo"
$wsReadme.Cells.Item(1, 5).Value = "Date"

# Helper: write a numeric-looking string ("20220302") as literal text
# without Excel auto-converting it to a number (which would also bump the
# cell's style). Routing it through a throwaway formula cell and reading
# the evaluated .Value back gives us a plain text value we can re-assign.
$scratch = $wsReadme.Cells.Item(200, 200)
$scratch.Formula = '="20220302"'
$newDate = $scratch.Value
$scratch.Clear()

$sheetNames = @(
    "Criteria Failing, Air Speed 0.1",
    "Criterion 1, Air Speed 0.1",
    "Criterion 3, Air Speed 0.1",
    "Criterion 2, Air Speed 0.1"
)

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $wsReadme.Cells.Item($r, 2).Value = "jovyan"
    $wsReadme.Cells.Item($r, 3).Value = $sheetNames[$i]
    $wsReadme.Cells.Item($r, 4).Value = "/c/e"
    $wsReadme.Cells.Item($r, 5).Value = $newDate
}
